# Applies the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.485.20"
$ws.Range("E2").Value = "  +0.71%  "

# Row 3
$ws.Range("D3").Value = "1.727.72"
$ws.Range("E3").Value = "  +0.80%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.80"
$ws.Range("E5").Value = "  +2.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4796"
$ws.Range("E7").Value = "  +2.36%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2687"
$ws.Range("E8").Value = "  +2.45%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06228"
$ws.Range("E9").Value = "  +0.30%  "

# Row 10
$ws.Range("D10").Value = "1.727.55"
$ws.Range("E10").Value = "  +1.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07143"
$ws.Range("E11").Value = "  +1.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.75"
$ws.Range("E12").Value = "  +3.92%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6213"
$ws.Range("E13").Value = "  +6.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.521"
$ws.Range("E14").Value = "  +2.70%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.19"
$ws.Range("E15").Value = "  +1.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"

# Row 17
$ws.Range("D17").Value = "26.497.17"
$ws.Range("E17").Value = "  +0.76%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006952"
$ws.Range("E19").Value = "  +2.20%  "

# Row 20
$ws.Range("E20").Value = "  +1.08%  "

# Row 21
$ws.Range("D21").Value = "1.950.13"
$ws.Range("E21").Value = "  +0.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.538"
$ws.Range("E22").Value = "  -0.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.970"
$ws.Range("E23").Value = "  +2.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.300"
$ws.Range("E24").Value = "  -0.43%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.42"
$ws.Range("E25").Value = "  +0.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.35"
$ws.Range("E26").Value = "  +1.34%  "

# Row 27
$ws.Range("E27").Value = "  +2.62%  "

# Row 28
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.78"
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.979"
$ws.Range("E30").Value = "  -1.23%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08023"
$ws.Range("E31").Value = "  +3.77%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.728"
$ws.Range("E32").Value = "  +1.68%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04570"
$ws.Range("E33").Value = "  +3.76%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.617"
$ws.Range("E34").Value = "  +0.16%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6370"
$ws.Range("E35").Value = "  +2.97%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9915"
$ws.Range("E36").Value = "  +2.67%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9377"
$ws.Range("E37").Value = "  +2.53%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.104"
$ws.Range("E38").Value = "  +11.13%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.415"
$ws.Range("E39").Value = "  +0.46%  "

# Row 40
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.15"
$ws.Range("E40").Value = "  -6.84%  "

# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("E41").Value = "  +0.51%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.707"
$ws.Range("E42").Value = "  +9.68%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01501"
$ws.Range("E43").Value = "  +2.63%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3917"
$ws.Range("E44").Value = "  +3.21%  "

# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.967"
$ws.Range("E45").Value = "  +11.86%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1191"
$ws.Range("E46").Value = "  +4.20%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05328"
$ws.Range("E47").Value = "  +0.68%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "31.04"
$ws.Range("E48").Value = "  +1.61%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.902"
$ws.Range("E49").Value = "  +2.96%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.266"
$ws.Range("E50").Value = "  +4.06%  "

# Row 51
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3433"
$ws.Range("E51").Value = "  +1.94%  "
